{"js": "// Resume education bullet originally reads:\n//   \"Working Towards a Bachelor of Science in Electrical Engineering\"\n// It needs to read:\n//   \"Working towards a bachelor of science in electrical engineering,\"\n// i.e. de-capitalize \"Towards\"/\"Bachelor\"/\"Science\"/\"Electrical\"/\"Engineering\"\n// down to \"towards\"/\"bachelor\"/\"science\"/\"electrical\"/\"engineering\" and add a\n// trailing comma (the line continues into the school name on the next line).\n\nconst body = context.document.body;\n\nconst phraseResults = body.search(\n  \"Working Towards a Bachelor of Science in Electrical Engineering\",\n  { matchCase: true, matchWholeWord: false }\n);\nphraseResults.load(\"items\");\nawait context.sync();\n\nif (phraseResults.items.length === 0) {\n  throw new Error(\"Could not find the degree line to update.\");\n}\n\nconst lineRange = phraseResults.items[0];\n\n// Lower-case the leading capital letter of each of these words, one at a\n// time (mirrors the series of small in-place corrections that produced the\n// change) using a search scoped to the line itself so we never touch a\n// same-letter match elsewhere in the document.\nconst wordsToLowerCase = [\"Towards\", \"Bachelor\", \"Science\", \"Electrical\", \"Engineering\"];\n\nfor (const word of wordsToLowerCase) {\n  const wordResults = lineRange.search(word, { matchCase: true });\n  wordResults.load(\"items\");\n  await context.sync();\n\n  if (wordResults.items.length === 0) {\n    throw new Error(`Could not find \"${word}\" on the degree line.`);\n  }\n\n  const wordRange = wordResults.items[0];\n  const firstLetterResults = wordRange.search(word.charAt(0), { matchCase: true });\n  firstLetterResults.load(\"items\");\n  await context.sync();\n\n  firstLetterResults.items[0].insertText(word.charAt(0).toLowerCase(), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Append the trailing comma at the very end of the (now lower-cased) line.\nlineRange.insertText(\",\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Resume education bullet originally reads:\n#   \"Working Towards a Bachelor of Science in Electrical Engineering\"\n# It needs to read:\n#   \"Working towards a bachelor of science in electrical engineering,\"\n# i.e. de-capitalize \"Towards\"/\"Bachelor\"/\"Science\"/\"Electrical\"/\"Engineering\"\n# down to \"towards\"/\"bachelor\"/\"science\"/\"electrical\"/\"engineering\" and add a\n# trailing comma (the line continues into the school name on the next line).\n\n$d = $word.ActiveDocument\n\n# Locate the degree line once; every later Find is scoped to a duplicate of\n# this range so it can't drift onto a same-letter match elsewhere in the\n# resume.\n$lineRange = $d.Content.Duplicate\n$lineFind = $lineRange.Find\n$lineFind.ClearFormatting()\n$lineFind.Text = \"Working Towards a Bachelor of Science in Electrical Engineering\"\n$lineFind.MatchCase = $true\n$lineFind.Forward = $true\n$lineFind.Wrap = 0\nif (-not $lineFind.Execute()) {\n    throw \"Could not find the degree line to update.\"\n}\n\n# Lower-case the leading capital letter of each of these words, one at a\n# time (mirrors the series of small in-place corrections that produced the\n# change). NOTE: the loop variable is deliberately NOT named `$word` so it\n# can't shadow the pre-seeded Word Application COM object.\n$targetWords = @(\"Towards\", \"Bachelor\", \"Science\", \"Electrical\", \"Engineering\")\nforeach ($targetWord in $targetWords) {\n    $wordRange = $lineRange.Duplicate\n    $wf = $wordRange.Find\n    $wf.ClearFormatting()\n    $wf.Text = $targetWord\n    $wf.MatchCase = $true\n    $wf.Forward = $true\n    $wf.Wrap = 0\n    if (-not $wf.Execute()) {\n        throw \"Could not find '$targetWord' on the degree line.\"\n    }\n\n    $letterRange = $wordRange.Duplicate\n    $lf = $letterRange.Find\n    $lf.ClearFormatting()\n    $lf.Text = $targetWord.Substring(0, 1)\n    $lf.MatchCase = $true\n    $lf.Forward = $true\n    $lf.Wrap = 0\n    if (-not $lf.Execute()) {\n        throw \"Could not find the leading letter of '$targetWord'.\"\n    }\n\n    $letterRange.Text = $targetWord.Substring(0, 1).ToLower()\n}\n\n# Append the trailing comma right after \"Engineering\" (still inside the\n# paragraph, before its end-of-paragraph mark). Re-find the line fresh so\n# the insertion point reflects the edits made above.\n$finalRange = $d.Content.Duplicate\n$finalFind = $finalRange.Find\n$finalFind.ClearFormatting()\n$finalFind.Text = \"Working towards a bachelor of science in electrical engineering\"\n$finalFind.MatchCase = $true\n$finalFind.Forward = $true\n$finalFind.Wrap = 0\nif (-not $finalFind.Execute()) {\n    throw \"Could not find the updated degree line to append the comma.\"\n}\n$finalRange.Collapse(0)\n$finalRange.InsertAfter(\",\")\n"}
